$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value looks like a plain number need to be forced
# to remain text (matching the original inlineStr/shared-string type) so Excel
# does not silently convert them to numeric cells.
$textCells = @('D5', 'D7', 'D10', 'D12', 'D13', 'D15', 'D16', 'D19', 'D21', 'D26', 'D28', 'D32', 'D34', 'D35', 'D36', 'D38', 'D40', 'D41', 'D44', 'D46', 'D50')
foreach ($cellref in $textCells) {
    $ws.Range($cellref).NumberFormat = "@"
}

$ws.Range('D2').Value = '36.507.76'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').Value = '1.954.45'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '243.18'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('E6').Value = '  +2.96%  '
$ws.Range('D7').Value = '60.28'
$ws.Range('E7').Value = '  +6.94%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  +5.27%  '
$ws.Range('D10').Value = '0.0789'
$ws.Range('E10').Value = '  -2.03%  '
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').Value = '14.15'
$ws.Range('E12').Value = '  +6.82%  '
$ws.Range('D13').Value = '0.839'
$ws.Range('E13').Value = '  +4.86%  '
$ws.Range('D14').Value = '2.241.70'
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('D15').Value = '21.50'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').Value = '5.27'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').Value = '1.952.52'
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('D18').Value = '36.464.90'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').Value = '69.23'
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('D21').Value = '229.29'
$ws.Range('E21').Value = '  +1.46%  '
$ws.Range('E22').Value = '  +3.18%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  +3.29%  '
$ws.Range('E25').Value = '  +4.00%  '
$ws.Range('D26').Value = '0.142'
$ws.Range('E26').Value = '  +8.20%  '
$ws.Range('E27').Value = '  +1.19%  '
$ws.Range('D28').Value = '160.53'
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range('E29').Value = '  +1.60%  '
$ws.Range('E30').Value = '  +21.71%  '
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('D32').Value = '4.76'
$ws.Range('E32').Value = '  +4.82%  '
$ws.Range('E33').Value = '  +0.52%  '
$ws.Range('D34').Value = '4.46'
$ws.Range('E34').Value = '  +8.28%  '
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '3.43'
$ws.Range('E36').Value = '  +7.21%  '
$ws.Range('E37').Value = '  +3.76%  '
$ws.Range('D38').Value = '1.77'
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('E39').Value = '  -10.62%  '
$ws.Range('D40').Value = '0.0966'
$ws.Range('E40').Value = '  -1.54%  '
$ws.Range('D41').Value = '2.91'
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('E42').Value = '  +2.57%  '
$ws.Range('E43').Value = '  +1.38%  '
$ws.Range('D44').Value = '15.86'
$ws.Range('E44').Value = '  +1.04%  '
$ws.Range('D45').Value = '1.361.82'
$ws.Range('E45').Value = '  +2.63%  '
$ws.Range('D46').Value = '88.73'
$ws.Range('E46').Value = '  +3.84%  '
$ws.Range('E47').Value = '  +0.68%  '
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('E49').Value = '  +0.98%  '
$ws.Range('D50').Value = '45.67'
$ws.Range('E50').Value = '  +6.62%  '
$ws.Range('D51').Value = '2.137.09'
